# Update cryptos list values to match the latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text interpretation so numeric-looking strings (e.g. "243.74")
    # are not silently coerced into floating point numbers by Excel,
    # matching the original inline-string cell content.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "36.224.78"
$ws.Cells.Item(2, 5).Value = "  -4.13%  "

$ws.Cells.Item(3, 4).Value = "1.970.45"
$ws.Cells.Item(3, 5).Value = "  -4.17%  "

$ws.Cells.Item(4, 5).Value = "  +0.01%  "

Set-TextCell 5 4 "243.74"
$ws.Cells.Item(5, 5).Value = "  -3.82%  "

Set-TextCell 6 4 "0.624"
$ws.Cells.Item(6, 5).Value = "  -4.06%  "

Set-TextCell 7 4 "59.75"
$ws.Cells.Item(7, 5).Value = "  -9.33%  "

$ws.Cells.Item(8, 5).Value = "  -0.03%  "

Set-TextCell 9 4 "0.376"
$ws.Cells.Item(9, 5).Value = "  -1.28%  "

Set-TextCell 10 4 "57.46"
$ws.Cells.Item(10, 5).Value = "  -5.48%  "

Set-TextCell 11 4 "0.0805"
$ws.Cells.Item(11, 5).Value = "  +5.09%  "

$ws.Cells.Item(12, 5).Value = "  -1.92%  "

$ws.Cells.Item(13, 2).Value = "Avalanche"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 13 4 "22.96"
$ws.Cells.Item(13, 5).Value = "  +10.19%  "

$ws.Cells.Item(14, 2).Value = "Polygon"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 14 4 "0.854"
$ws.Cells.Item(14, 5).Value = "  -7.31%  "

$ws.Cells.Item(15, 2).Value = "Chainlink"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 15 4 "14.06"
$ws.Cells.Item(15, 5).Value = "  -6.01%  "

$ws.Cells.Item(16, 4).Value = "2.257.21"
$ws.Cells.Item(16, 5).Value = "  -4.24%  "

Set-TextCell 17 4 "5.41"
$ws.Cells.Item(17, 5).Value = "  -3.67%  "

$ws.Cells.Item(18, 4).Value = "1.968.31"
$ws.Cells.Item(18, 5).Value = "  -4.38%  "

$ws.Cells.Item(19, 4).Value = "36.157.08"
$ws.Cells.Item(19, 5).Value = "  -4.11%  "

Set-TextCell 20 4 "71.19"
$ws.Cells.Item(20, 5).Value = "  -3.88%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0858"
$ws.Cells.Item(21, 5).Value = "  -2.38%  "

$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 22 4 "5.25"
$ws.Cells.Item(22, 5).Value = "  -3.23%  "

$ws.Cells.Item(23, 2).Value = "BitcoinCash"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 23 4 "233.55"
$ws.Cells.Item(23, 5).Value = "  -3.06%  "

$ws.Cells.Item(24, 5).Value = "  +0.07%  "

Set-TextCell 25 4 "2.58"
$ws.Cells.Item(25, 5).Value = "  -3.18%  "

Set-TextCell 26 4 "2.29"
$ws.Cells.Item(26, 5).Value = "  -4.33%  "

Set-TextCell 27 4 "9.77"
$ws.Cells.Item(27, 5).Value = "  +1.16%  "

Set-TextCell 28 4 "161.33"
$ws.Cells.Item(28, 5).Value = "  -0.69%  "

$ws.Cells.Item(29, 2).Value = "Kaspa"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell 29 4 "0.135"
$ws.Cells.Item(29, 5).Value = "  +15.22%  "

$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 30 4 "19.76"
$ws.Cells.Item(30, 5).Value = "  -1.55%  "

Set-TextCell 31 4 "0.120"
$ws.Cells.Item(31, 5).Value = "  -1.66%  "

Set-TextCell 32 4 "4.88"
$ws.Cells.Item(32, 5).Value = "  -8.56%  "

$ws.Cells.Item(33, 5).Value = "  -8.27%  "

Set-TextCell 34 4 "0.0620"
$ws.Cells.Item(34, 5).Value = "  -0.94%  "

Set-TextCell 35 4 "4.46"
$ws.Cells.Item(35, 5).Value = "  -6.00%  "

$ws.Cells.Item(36, 5).Value = "  +0.11%  "

Set-TextCell 37 4 "2.27"
$ws.Cells.Item(37, 5).Value = "  -7.37%  "

$ws.Cells.Item(38, 2).Value = "THORChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell 38 4 "6.03"
$ws.Cells.Item(38, 5).Value = "  -1.87%  "

$ws.Cells.Item(39, 2).Value = "WEMIXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell 39 4 "1.80"
$ws.Cells.Item(39, 5).Value = "  -2.72%  "

Set-TextCell 40 4 "3.08"
$ws.Cells.Item(40, 5).Value = "  +7.23%  "

$ws.Cells.Item(41, 2).Value = "TrustWalletToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 41 4 "1.23"
$ws.Cells.Item(41, 5).Value = "  -1.32%  "

$ws.Cells.Item(42, 2).Value = "Cronos"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell 42 4 "0.0977"
$ws.Cells.Item(42, 5).Value = "  -8.41%  "

Set-TextCell 43 4 "2.89"
$ws.Cells.Item(43, 5).Value = "  -0.57%  "

Set-TextCell 44 4 "0.0213"
$ws.Cells.Item(44, 5).Value = "  -3.68%  "

Set-TextCell 45 4 "1.09"
$ws.Cells.Item(45, 5).Value = "  -5.37%  "

Set-TextCell 46 4 "92.12"
$ws.Cells.Item(46, 5).Value = "  -3.68%  "

$ws.Cells.Item(47, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 47 4 "16.14"
$ws.Cells.Item(47, 5).Value = "  -6.51%  "

$ws.Cells.Item(48, 2).Value = "FraxShare"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell 48 4 "7.57"
$ws.Cells.Item(48, 5).Value = "  -4.71%  "

$ws.Cells.Item(49, 4).Value = "1.350.84"
$ws.Cells.Item(49, 5).Value = "  -4.14%  "

Set-TextCell 50 4 "2.84"
$ws.Cells.Item(50, 5).Value = "  -3.72%  "

$ws.Cells.Item(51, 2).Value = "MultiversX"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextCell 51 4 "44.33"
$ws.Cells.Item(51, 5).Value = "  -6.01%  "
